$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.916.46'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.31%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.495.54'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.07%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.39%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.62%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.592'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.16%  '

# Row 9
$ws.Range("E9").Value = '  +3.23%  '

# Row 10
$ws.Range("E10").Value = '  -1.17%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.433'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.30%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.098.75'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.24%  '

# Row 13
$ws.Range("E13").Value = '  -0.26%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.20%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '66.932.80'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.35%  '

# Row 16
$ws.Range("E16").Value = '  +0.12%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.474.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.34%  '

# Row 18
$ws.Range("E18").Value = '  -0.66%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.95%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '395.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.60%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.97'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.36%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.32%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.06%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.536'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.01%  '

# Row 25
$ws.Range("E25").Value = '  -0.81%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.12%  '

# Row 27
$ws.Range("E27").Value = '  -0.29%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.17%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.86%  '

# Row 30
$ws.Range("E30").Value = '  -2.43%  '

# Row 31
$ws.Range("E31").Value = '  -0.53%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.74'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.71%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.36'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.76%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.63'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.31%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '162.81'
$ws.Range("D35").Style = "Normal"

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.878'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.89%  '

# Row 37
$ws.Range("E37").Value = '  -1.16%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.93'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.61%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.66'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.09%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0741'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.45%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '27.27'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.45%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.832.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.35%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.28'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.37%  '

# Row 44
$ws.Range("E44").Value = '  -0.99%  '

# Row 45
$ws.Range("E45").Value = '  +2.30%  '

# Row 46
$ws.Range("E46").Value = '  -3.67%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '337.49'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.68%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '34.72'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.02%  '

# Row 49
$ws.Range("E49").Value = '  -1.54%  '

# Row 50
$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.842'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.24%  '

# Row 51
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.42'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.90%  '
